$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "value" header (B1) to "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# 2. Replace the single data row (row 2) with the full yearly series.
#    Dates are year-end (Dec 31) values from 2005 through 2025; the first
#    and last rows only carry a date, with no corresponding value.
$dates = @(
    "2005-12-31", "2006-12-31", "2007-12-31", "2008-12-31", "2009-12-31",
    "2010-12-31", "2011-12-31", "2012-12-31", "2013-12-31", "2014-12-31",
    "2015-12-31", "2016-12-31", "2017-12-31", "2018-12-31", "2019-12-31",
    "2020-12-31", "2021-12-31", "2022-12-31", "2023-12-31", "2024-12-31",
    "2025-12-31"
)

$values = @(
    $null,
    12.85455285386146,
    8.323723251380377,
    2.675214973655216,
    -14.07263945793084,
    14.00662378688902,
    10.01426242069761,
    5.166028195387984,
    0.3163778774614823,
    4.811464743291949,
    4.153025533745458,
    2.211964547984113,
    5.059349743581909,
    2.287635922746656,
    0.6462611928503614,
    -12.47081270006417,
    5.490291529373104,
    0.01028356335206482,
    -0.7275558254695946,
    0.2831497518338555,
    $null
)

# Copy the existing date cell's formatting (A2) first, so every new date
# cell down through row 22 keeps the same number format/font/border as the
# original template row, then fill in the values.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2:A22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]

    if ($null -ne $values[$i]) {
        $ws.Cells.Item($row, 2).Value = $values[$i]
    } else {
        $ws.Cells.Item($row, 2).ClearContents() | Out-Null
    }
}
